$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F, matching the style used by the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Add time_taken values for each data row (plain, unstyled cells like the rest of column F)
$times = @(
    "2021-10-05 10:50:19.045354",
    "2021-10-05 10:50:19.045366",
    "2021-10-05 10:50:19.045370",
    "2021-10-05 10:50:19.045374",
    "2021-10-05 10:50:19.045377",
    "2021-10-05 10:50:19.045381",
    "2021-10-05 10:50:19.045384",
    "2021-10-05 10:50:19.045387",
    "2021-10-05 10:50:19.045390",
    "2021-10-05 10:50:19.045394",
    "2021-10-05 10:50:19.045397"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
